$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Extend the repeating 4-column "Assignment/Grade/Comments/gap"
#    block (currently ending at column Y:AB) seven more times so the
#    sheet covers columns A:BD (matching the grading-sheet template
#    used for every assignment column group).
# ------------------------------------------------------------------
$srcBlock = $ws.Range("Y1:AB11")
for ($i = 0; $i -lt 7; $i++) {
    $destCell = $ws.Cells.Item(1, 29 + $i * 4)
    $srcBlock.Copy($destCell)
}

# ------------------------------------------------------------------
# 2) Fill in the grades for assignments "3 CPP", "4 UE" and "5 UE"
#    for the student (row 2). Assignment 6 (Z:AB) is left ungraded.
# ------------------------------------------------------------------
$ws.Range("N2").Value = "3 CPP"
$ws.Range("O2").Value = 80
$ws.Range("P2").Value = "Good!, see my comments"

$ws.Range("R2").Value = "4 UE"
$ws.Range("S2").Value = 100
$ws.Range("T2").Value = "Excellent!"

$ws.Range("V2").Value = "5 UE"
$ws.Range("W2").Value = 33
$ws.Range("X2").Value = "did not create another actor type like Arrow and Target (see items 2 and 3 in class 5 HW)"

# ------------------------------------------------------------------
# 3) Widen the "5 UE" comments column (X) so the long comment is
#    readable.
# ------------------------------------------------------------------
$ws.Columns.Item(24).ColumnWidth = 20.6

# ------------------------------------------------------------------
# 4) Taller header / data rows so the wrapped assignment headers and
#    the long new comment are fully visible.
# ------------------------------------------------------------------
$ws.Rows.Item(1).RowHeight = 30
$ws.Rows.Item(2).RowHeight = 60

# ------------------------------------------------------------------
# 5) Leave the selection where the editor ended up.
# ------------------------------------------------------------------
$ws.Range("AS18").Select()
